$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "####1er Torneo Federativo - C.A.E. - Sub 23, Prejuveniles y sub 23 (28 de Febrero y 1 de Marzo) - Juniors (Domingo 1 de Marzo)"
$ws.Range("B3").Value = "Prejuveniles"
$ws.Range("C3").Value = "caballeros"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "Petric, Juan Cruz"
$ws.Range("F3").Value = 82
$ws.Range("G3").Value = "OK"
